$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 17: replace CyVerse / Tucson / Arizona entry with ESIIL / Boulder / Colorado entry
# (set B/C before A so new shared-string entries are appended in ESIIL, Boulder, full-name order)
$ws.Range("B17").Value = "ESIIL"
$ws.Range("C17").Value = "Boulder"
$ws.Range("A17").Value = "Environmental Data Science Innovation & Inclusion Lab"
$ws.Range("D17").Value = "Colorado"
$ws.Range("E17").Value = "United States of America"
$ws.Range("F17").Value = "USA"
$ws.Range("G17").Value = 40.014986
$ws.Range("H17").Value = -105.270546
$ws.Range("I17").Value = "Yes"

# Update the selected cell to reflect the active cell used when saving
$ws.Range("H17").Select()
